$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 154: headers for the new "test runs" batch table ----
$ws.Cells.Item(154, 1).Value = "Gecode"
$ws.Cells.Item(154, 2).Value = "16th Jun"
$ws.Cells.Item(154, 3).Value = "5 min"

# ---- Row 155: Matrix size ----
$ws.Cells.Item(155, 1).Value = "Matrix size"
$matrixSizes = @(104,204,511,1037,1568,2149,2657,3075,3505,4025,5043,6085)
for ($i = 0; $i -lt $matrixSizes.Length; $i++) {
    $ws.Cells.Item(155, 2 + $i).Value = $matrixSizes[$i]
}

# (finish row 154 header row after "Matrix size" is registered, to match shared-string order)
$ws.Cells.Item(154, 4).Value = "test runs"

# ---- Row 156: cost £ ----
$ws.Cells.Item(156, 1).Value = "cost £"
$costVals = @(851,1821,4617,9044,14108,18991,23556,27895,31368)
for ($i = 0; $i -lt $costVals.Length; $i++) {
    $ws.Cells.Item(156, 2 + $i).Value = $costVals[$i]
}
$ws.Cells.Item(156, 11).Value = "unable to compute"
$ws.Cells.Item(156, 12).Value = "unable to compute"
$ws.Cells.Item(156, 13).Value = "unable to compute"

# ---- Row 157: emissions kg ----
$ws.Cells.Item(157, 1).Value = "emissions kg"
$emissionVals = @(865,1646,4454,8634,12212,17782,21449,25804,28934)
for ($i = 0; $i -lt $emissionVals.Length; $i++) {
    $ws.Cells.Item(157, 2 + $i).Value = $emissionVals[$i]
}
$ws.Cells.Item(157, 11).Value = "unable to compute"
$ws.Cells.Item(157, 12).Value = "unable to compute"
$ws.Cells.Item(157, 13).Value = "unable to compute"

# ---- Row 158: food waste sum nutrients ----
$ws.Cells.Item(158, 1).Value = "food waste sum nutrients"
$wasteVals = @(58649,79523,182057,401520,617231,763670,1104941,1039116,1490343)
for ($i = 0; $i -lt $wasteVals.Length; $i++) {
    $ws.Cells.Item(158, 2 + $i).Value = $wasteVals[$i]
}
$ws.Cells.Item(158, 11).Value = "unable to compute"
$ws.Cells.Item(158, 12).Value = "unable to compute"
$ws.Cells.Item(158, 13).Value = "unable to compute"

# ---- Row 159: cost/(days x people) ----
$ws.Cells.Item(159, 1).Value = "cost/(days x people)"
$costFormulas = @("=851/104","=1821/204","=4617/511","=9044/1037","=14108/1568","=18991/2149","=23556/2657","=27895/3075","=31368/3505")
for ($i = 0; $i -lt $costFormulas.Length; $i++) {
    $ws.Cells.Item(159, 2 + $i).Formula = $costFormulas[$i]
}
$ws.Cells.Item(159, 11).Value = "unable to compute"
$ws.Cells.Item(159, 12).Value = "unable to compute"
$ws.Cells.Item(159, 13).Value = "unable to compute"

# ---- Row 160: emissions /(days x people) ----
$ws.Cells.Item(160, 1).Value = "emissions /(days x people)"
$emissionFormulas = @("=865/104","=1646/204","=4454/511","=8634/1037","=12212/1568","=17782/2149","=21449/2657","=25804/3075","=28934/3505")
for ($i = 0; $i -lt $emissionFormulas.Length; $i++) {
    $ws.Cells.Item(160, 2 + $i).Formula = $emissionFormulas[$i]
}
$ws.Cells.Item(160, 11).Value = "unable to compute"
$ws.Cells.Item(160, 12).Value = "unable to compute"
$ws.Cells.Item(160, 13).Value = "unable to compute"

# ---- Row 161: food waste /(days x people) ----
$ws.Cells.Item(161, 1).Value = "food waste /(days x people)"
$wasteFormulas = @("=58649/104","=79523/204","=182057/511","=401520/1037","=617231/1568","=763670/2149","=1104941/2657","=1039116/3075","=1490343/3505")
for ($i = 0; $i -lt $wasteFormulas.Length; $i++) {
    $ws.Cells.Item(161, 2 + $i).Formula = $wasteFormulas[$i]
}
$ws.Cells.Item(161, 11).Value = "unable to compute"
$ws.Cells.Item(161, 12).Value = "unable to compute"
$ws.Cells.Item(161, 13).Value = "unable to compute"

# Row 162 intentionally left blank (gap before actual run summary)

# ---- Row 163: actual run ----
$ws.Cells.Item(163, 1).Value = "actual run"

# ---- Row 164: Matrix size (actual) ----
$ws.Cells.Item(164, 1).Value = "Matrix size"
$ws.Cells.Item(164, 2).Value = 3505

# ---- Row 165: cost £ (actual) ----
$ws.Cells.Item(165, 1).Value = "cost £"
$ws.Cells.Item(165, 2).Value = 25848

# ---- Row 166: emissions kg (actual) ----
$ws.Cells.Item(166, 1).Value = "emissions kg"
$ws.Cells.Item(166, 2).Value = 22713

# ---- Row 167: food waste sum nutrients (actual) ----
$ws.Cells.Item(167, 1).Value = "food waste sum nutrients"
$ws.Cells.Item(167, 2).Value = 776138

# ---- Row 168: cost/(days x people) (actual) ----
$ws.Cells.Item(168, 1).Value = "cost/(days x people)"
$ws.Cells.Item(168, 2).Formula = "=25848/3505"

# ---- Row 169: emissions /(days x people) (actual) ----
$ws.Cells.Item(169, 1).Value = "emissions /(days x people)"
$ws.Cells.Item(169, 2).Formula = "=22713/3505"

# ---- Row 170: food waste /(days x people) (actual) ----
$ws.Cells.Item(170, 1).Value = "food waste /(days x people)"
$ws.Cells.Item(170, 2).Formula = "=776138/3505"

# ---- Update sheet view to match the new scroll/selection position ----
$ws.Application.ActiveWindow.ScrollRow = 153
$ws.Range("F164").Select()
